# Nipasha_discountingtask.xlsx — revert a prior sort on A2:E8 back to the
# original (pre-sort) row order, tidy up the D-column labels and the
# "k" header text, rebuild the G/H helper formulas on row 5 instead of
# row 6, and restore the previous selection / column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Restore the original (unsorted) row order for the data block
#    A2:D8, and rewrite the E-column ratio formula.
# ---------------------------------------------------------------------
$data = @(
    @(15, 35, 13, "l"),
    @(67, 75, 119, "t"),
    @(27, 50, 21, "t"),
    @(11, 30, 7, "l"),
    @(40, 55, 62, "t"),
    @(78, 80, 162, "t"),
    @(49, 60, 89, "t")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value2 = $vals[0]
    $ws.Cells.Item($row, 2).Value2 = $vals[1]
    $ws.Cells.Item($row, 3).Value2 = $vals[2]
    $ws.Cells.Item($row, 4).Value2 = $vals[3]
}

# Row 2's ratio formula stays a standalone (non-shared) formula.
$ws.Range("E2").Formula = "=((B2/A2)-1)/C2"

# Rows 3-8 become one shared-formula group.
$ws.Range("E3:E8").Formula = "=((B3/A3)-1)/C3"

# ---------------------------------------------------------------------
# 2. Rebuild the GEOMEAN / SQRT helper cells on row 5 (they used to
#    live on row 6, next to the old sorted position of this data).
# ---------------------------------------------------------------------
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("G5").Formula = "=GEOMEAN(E5,E6)"
$ws.Range("H5").FormulaArray = "=SQRT(E5:E6)"

# ---------------------------------------------------------------------
# 3. Tidy up shared-string text: "k (discounting parameter)" -> "k".
# ---------------------------------------------------------------------
$ws.Range("E1").Value2 = "k"

# ---------------------------------------------------------------------
# 4. Drop the leftover sort state left behind on the sheet.
# ---------------------------------------------------------------------
$ws.Sort.SortFields.Clear()

# ---------------------------------------------------------------------
# 5. Column width tweaks.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 10.736979166666666
$ws.Columns.Item(5).ColumnWidth = 13.02213541666666

# ---------------------------------------------------------------------
# 6. Restore the previous cell selection.
# ---------------------------------------------------------------------
$ws.Range("D10").Select() | Out-Null
